# Edit script:
#  1. Re-point the three table styles (slides 14/15/16) from the
#     "Table_0" style {B2D7E0AC-39CF-4BEF-B9AE-8AFDB70F7FAF} to
#     {7FDC3200-D55B-46D8-BDFB-1EA3F7845142}.
#  2. Re-colour the (only reachable) theme - ppt/theme/theme1.xml,
#     used by the slide master - from the "Integral / Red Violet"
#     colour scheme to the "Office" colour scheme, matching the
#     target theme swap.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{7FDC3200-D55B-46D8-BDFB-1EA3F7845142}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
function Set-SchemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$master = $p.SlideMaster
$cs = $master.ColorScheme

Set-SchemeColor $cs 1  "000000"  # dk1
Set-SchemeColor $cs 2  "FFFFFF"  # lt1
Set-SchemeColor $cs 3  "44546A"  # dk2
Set-SchemeColor $cs 4  "E7E6E6"  # lt2
Set-SchemeColor $cs 5  "5B9BD5"  # accent1
Set-SchemeColor $cs 6  "ED7D31"  # accent2
Set-SchemeColor $cs 7  "A5A5A5"  # accent3
Set-SchemeColor $cs 8  "FFC000"  # accent4
Set-SchemeColor $cs 9  "4472C4"  # accent5
Set-SchemeColor $cs 10 "70AD47"  # accent6
Set-SchemeColor $cs 11 "0563C1"  # hlink
Set-SchemeColor $cs 12 "954F72"  # folHlink
